$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.798.38'
$ws.Range("E2").Value = '  +0.80%  '
$ws.Range("E3").Value = '  +0.76%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '303.92'
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.61'
$ws.Range("E6").Value = '  +0.91%  '
$ws.Range("E7").Value = '  +1.93%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.484'
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.57'
$ws.Range("E10").Value = '  +1.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.72'
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("E12").Value = '  +0.36%  '
$ws.Range("E13").Value = '  -1.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.68'
$ws.Range("E14").Value = '  +1.20%  '
$ws.Range("D15").Value = '2.623.66'
$ws.Range("E15").Value = '  +0.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.27'
$ws.Range("E16").Value = '  +0.80%  '
$ws.Range("D17").Value = '2.271.27'
$ws.Range("E17").Value = '  -6.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.776'
$ws.Range("E18").Value = '  +3.32%  '
$ws.Range("D19").Value = '41.718.62'
$ws.Range("E19").Value = '  +0.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.52'
$ws.Range("E20").Value = '  +2.76%  '
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.12'
$ws.Range("E23").Value = '  +0.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '244.00'
$ws.Range("E24").Value = '  +1.30%  '
$ws.Range("E25").Value = '  +0.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.93'
$ws.Range("E26").Value = '  +3.90%  '
$ws.Range("E27").Value = '  +0.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.97'
$ws.Range("E28").Value = '  +0.93%  '
$ws.Range("E29").Value = '  -1.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.34'
$ws.Range("E31").Value = '  +4.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '160.49'
$ws.Range("E32").Value = '  +1.51%  '
$ws.Range("E33").Value = '  +1.10%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0742'
$ws.Range("E35").Value = '  +0.72%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.00'
$ws.Range("E36").Value = '  -1.27%  '
$ws.Range("E37").Value = '  +1.60%  '
$ws.Range("E38").Value = '  -0.20%  '
$ws.Range("E39").Value = '  +1.55%  '
$ws.Range("E40").Value = '  +0.82%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.90'
$ws.Range("E42").Value = '  -1.92%  '
$ws.Range("D43").Value = '2.000.71'
$ws.Range("E43").Value = '  -3.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.54'
$ws.Range("E44").Value = '  -4.20%  '
$ws.Range("E45").Value = '  +1.57%  '
$ws.Range("E46").Value = '  +1.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.14'
$ws.Range("E47").Value = '  +2.69%  '
$ws.Range("E48").Value = '  -2.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '52.64'
$ws.Range("E50").Value = '  +0.69%  '
$ws.Range("E51").Value = '  -1.09%  '
